$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that look numeric/date need to be forced to Text so Excel
# doesn't silently convert them to a number / date serial (matches the
# source data which stores everything as literal text).
$ws.Range("F19").NumberFormat = "@"
$ws.Range("I19").NumberFormat = "@"
$ws.Range("J19").NumberFormat = "@"
$ws.Range("K19").NumberFormat = "@"

$ws.Range("A19").Value = "Aydinova Narine Sergeevna"
$ws.Range("B19").Value = "Yurisprudensiya"
$ws.Range("C19").Value = "Rus tili"
$ws.Range("D19").Value = "Kunduzgi"
$ws.Range("E19").Value = "AD6055389"
$ws.Range("F19").Value = "60402085220078"
$ws.Range("G19").Value = "Toshkent shahri"
$ws.Range("H19").Value = "Mirzo Ulugʻbek tumani"
$ws.Range("I19").Value = "998909340132"
$ws.Range("J19").Value = "+998935617938"
$ws.Range("K19").Value = "2025-04-25"
